$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "Contact" row (row 11); this shifts all subsequent
# rows up by one and shrinks the used range from A1:B22 to A1:B21.
$ws.Rows("11").Delete()

# Update Version value: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Update Date value
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive now has a value of "true" (stored as literal text, not a
# boolean) -- build it via a formula + paste-special-values so Excel doesn't
# auto-coerce the literal "true" into a Boolean TRUE cell.
$ws.Range("B14").Formula = '="true"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
